$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 (AE3:AH3) from "A" to "O"
$ws.Range("AE3:AH3").Value = "O"

# Update AI3 value (entry date) from 20000101 to 20150101
$ws.Range("AI3").Value = 20150101

# Delete row 4 entirely (it becomes blank / removed)
$ws.Rows("4:4").Delete()

# Update the selected cell / view to E10 (also resets any scrolled topLeftCell)
$ws.Range("E10").Select()
